# Applies the cryptos.xlsx symbol-list refresh: updates the Price (column D)
# and Volume(1h) (column E) text values for the affected rows to match the
# latest scrape. Cells store plain numeric-looking / percentage-looking
# strings as literal text (not real numbers), so each value is written with
# a leading apostrophe to force Excel to keep it as text instead of
# auto-converting it to a Number/Percentage cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'335.92"
$ws.Range("D3").Value = "'43.89"
$ws.Range("E3").Value = "'6.66%"
$ws.Range("E4").Value = "'2.77%"
$ws.Range("D5").Value = "'0.08331"
$ws.Range("E5").Value = "'1.40%"
$ws.Range("D6").Value = "'8.854"
$ws.Range("E6").Value = "'1.15%"
$ws.Range("D7").Value = "'4.522"
$ws.Range("E7").Value = "'0.61%"
$ws.Range("D8").Value = "'1.963"
$ws.Range("E8").Value = "'-1.86%"
$ws.Range("E9").Value = "'-1.89%"
$ws.Range("D10").Value = "'0.9459"
$ws.Range("E10").Value = "'2.67%"
$ws.Range("D11").Value = "'0.1243"
$ws.Range("E11").Value = "'-3.07%"
$ws.Range("D12").Value = "'0.1976"
$ws.Range("E12").Value = "'1.15%"
$ws.Range("D13").Value = "'0.09762"
$ws.Range("E13").Value = "'4.11%"
$ws.Range("D14").Value = "'0.04569"
$ws.Range("E14").Value = "'17.16%"
$ws.Range("E15").Value = "'0.83%"
$ws.Range("E16").Value = "'-0.62%"
$ws.Range("D17").Value = "'0.006032"
$ws.Range("E17").Value = "'-3.03%"
$ws.Range("D18").Value = "'3.499"
$ws.Range("E18").Value = "'1.51%"
$ws.Range("D19").Value = "'0.3506"
$ws.Range("E19").Value = "'0.59%"
$ws.Range("D20").Value = "'8.744"
$ws.Range("E20").Value = "'6.03%"
$ws.Range("D21").Value = "'0.1364"
$ws.Range("E21").Value = "'-0.64%"
$ws.Range("D23").Value = "'0.04426"
$ws.Range("E23").Value = "'0.67%"
$ws.Range("D24").Value = "'0.001263"
$ws.Range("E24").Value = "'0.52%"
$ws.Range("D25").Value = "'0.004352"
$ws.Range("E25").Value = "'1.07%"
$ws.Range("E26").Value = "'5.12%"
$ws.Range("D27").Value = "'0.0003998"
$ws.Range("E39").Value = "'0.32%"
$ws.Range("D40").Value = "'0.05768"
$ws.Range("E40").Value = "'6.87%"
$ws.Range("D41").Value = "'0.007937"
$ws.Range("E41").Value = "'1.67%"
$ws.Range("D43").Value = "'0.008981"
$ws.Range("E43").Value = "'0.44%"
$ws.Range("D44").Value = "'0.002173"
$ws.Range("E44").Value = "'0.09%"
$ws.Range("D45").Value = "'0.01042"
$ws.Range("E45").Value = "'-10.09%"
$ws.Range("D46").Value = "'0.00007297"
$ws.Range("E46").Value = "'8.15%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("D48").Value = "'0.003193"
$ws.Range("E48").Value = "'0.02%"
$ws.Range("D49").Value = "'0.002275"
$ws.Range("E49").Value = "'-0.26%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.13%"
